# Implements "dtype check" column: insert a new boolean column
# ("boolean_perfect") right before the existing "all_numbers" column
# (G) on the test_converters sheet, pushing the old G/H/I
# (all_numbers / ints_with_missing / float) columns one to the right
# (H/I/J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh column at G - this shifts former columns G,H,I (and
# their column-width overrides) to H,I,J automatically, and grows the
# sheet dimension/row spans to column J for us.
$ws.Columns.Item(7).Insert()

# The freshly inserted column inherits the formatting (and number
# format / style) of the column to its left (F). Strip that so the
# new boolean values below aren't stamped with an unrelated style.
$ws.Columns.Item(7).ClearFormats()
$ws.Columns.Item(7).ColumnWidth = $ws.Columns.Item(6).ColumnWidth()

# Header
$ws.Range("G1").Value = "boolean_perfect"

# Data - True for the "perfect" boolean column(s)/dtype rows, False
# once the data stops being a clean boolean.
$ws.Range("G2").Value = $true
$ws.Range("G3").Value = $true
$ws.Range("G4").Value = $true
$ws.Range("G5").Value = $true
$ws.Range("G6").Value = $false
$ws.Range("G7").Value = $false
$ws.Range("G8").Value = $false
$ws.Range("G9").Value = $false
$ws.Range("G10").Value = $false

# Match the selection left behind in the saved workbook.
$ws.Range("G7:G10").Select() | Out-Null
